# Actualizacion automatica del tracker
# Rellena resultado/profit de partidos ya resueltos y agrega los partidos nuevos.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resultados de partidos pendientes ---
$ws.Range("G6").Value = "Acierto"
$ws.Range("H6").Value = 0.83

$ws.Range("G7").Value = "Acierto"
$ws.Range("H7").Value = 2

$ws.Range("G9").Value = "Fallo"
$ws.Range("H9").Value = -1

$ws.Range("G11").Value = "Fallo"
$ws.Range("H11").Value = -1

$ws.Range("G13").Value = "Fallo"
$ws.Range("H13").Value = -1

# --- Nuevos partidos agregados al tracker ---
$ws.Range("A36").Value = 14633182
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "2025-09-08"
$ws.Range("C36").Value = "Vadym Ursu"
$ws.Range("D36").Value = "Clement Chidekh"
$ws.Range("E36").Value = "Gana Vadym Ursu"
$ws.Range("F36").Value = 3.75

$ws.Range("A37").Value = 14638045
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "2025-09-08"
$ws.Range("C37").Value = "Nikolay Vylegzhanin"
$ws.Range("D37").Value = "Maxence Bertimon"
$ws.Range("E37").Value = "Gana Maxence Bertimon"
$ws.Range("F37").Value = 2.2
